# The edit renames the "friendly name" metadata of three inline
# pictures that live in the document's headers/footers:
#   - footer (first page)   : image2.png -> image1.png   (Pearson logo)
#   - footer (default)      : image2.png -> image1.png   (Pearson logo)
#   - header (first page)   : image1.jpg -> image2.jpg   (BTEC logo)
#
# The picture's display Name is stored in two places inside each
# drawing: <wp:docPr name="..."/> and <pic:cNvPr name="..."/>. The
# high level InlineShape object only exposes a setter that reaches the
# first of those two (wp:docPr), so instead we round-trip the whole
# package through Document.WordOpenXML (flat-OPC XML containing every
# part's raw markup) and do a precise, scoped text substitution on the
# `name="..."` attribute of the picture shapes, leaving the actual
# image relationships (r:embed) and file names untouched.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

# Swap the two logo shape names with each other.
$xml = $xml.Replace('name="image2.png"', 'name="image1.png"')
$xml = $xml.Replace('name="image1.jpg"', 'name="image2.jpg"')

$d.WordOpenXML = $xml

Write-Output "renamed inline shapes"
